$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 188 (shifts existing rows 188:301 down to 189:302)
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with the new data record
$ws.Cells.Item(188, 1).Value = 8
$ws.Cells.Item(188, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(188, 3).Value = "Coquimbo"
$ws.Cells.Item(188, 4).Value = 45176
$ws.Cells.Item(188, 5).Value = 4
$ws.Cells.Item(188, 6).Value = 100112001
$ws.Cells.Item(188, 7).Value = "Berenjena"
$ws.Cells.Item(188, 8).Value = "Sin especificar"
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 460
$ws.Cells.Item(188, 11).Value = 9000
$ws.Cells.Item(188, 12).Value = 10000
$ws.Cells.Item(188, 13).Value = 9500
$ws.Cells.Item(188, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(188, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(188, 16).Value = 190
$ws.Cells.Item(188, 17).Value = 50
$ws.Cells.Item(188, 18).Value = "Hortaliza"
